# Refresh recomputed NATMI TPM-derived statistics for Reln -> Itga3 (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.041827
$ws.Range("H2").Value = 0.125481
$ws.Range("I2").Value = 0.006279874897961605
$ws.Range("J2").Value = 0.006279874897961606
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.731629
$ws.Range("N2").Value = 8.194887
$ws.Range("O2").Value = 0.5547800938501829
$ws.Range("P2").Value = 0.554780093850183
$ws.Range("Q2").Value = 0.114255846183
$ws.Range("R2").Value = 1.028302615647
$ws.Range("S2").Value = 0.003483949585258547
$ws.Range("T2").Value = 0.003483949585258548

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.041827
$ws.Range("H3").Value = 0.125481
$ws.Range("I3").Value = 0.006279874897961605
$ws.Range("J3").Value = 0.006279874897961606
$ws.Range("M3").Value = 0.06813733333333333
$ws.Range("O3").Value = 0.01383834927121065
$ws.Range("P3").Value = 0.01383834927121065
$ws.Range("Q3").Value = 0.002849980241333333
$ws.Range("R3").Value = 0.025649822172
$ws.Range("S3").Value = [double]"8.6903102217501E-05"
$ws.Range("T3").Value = [double]"8.690310221750103E-05"

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.041827
$ws.Range("H4").Value = 0.125481
$ws.Range("I4").Value = 0.006279874897961605
$ws.Range("J4").Value = 0.006279874897961606
$ws.Range("M4").Value = 2.124038666666666
$ws.Range("N4").Value = 6.372115999999999
$ws.Range("O4").Value = 0.4313815568786064
$ws.Range("P4").Value = 0.4313815568786064
$ws.Range("Q4").Value = 0.08884216531066666
$ws.Range("R4").Value = 0.799579487796
$ws.Range("S4").Value = 0.002709022210485557
$ws.Range("T4").Value = 0.002709022210485557

# Row 5
$ws.Range("I5").Value = 0.1693441751896972
$ws.Range("J5").Value = 0.1693441751896972
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.731629
$ws.Range("N5").Value = 8.194887
$ws.Range("O5").Value = 0.5547800938501829
$ws.Range("P5").Value = 0.554780093850183
$ws.Range("Q5").Value = 3.081042591905999
$ws.Range("R5").Value = 27.729383327154
$ws.Range("S5").Value = 0.09394877740472202
$ws.Range("T5").Value = 0.09394877740472206

# Row 6
$ws.Range("I6").Value = 0.1693441751896972
$ws.Range("J6").Value = 0.1693441751896972
$ws.Range("M6").Value = 0.06813733333333333
$ws.Range("O6").Value = 0.01383834927121065
$ws.Range("P6").Value = 0.01383834927121065
$ws.Range("Q6").Value = 0.07685305218933332
$ws.Range("R6").Value = 0.6916774697039999
$ws.Range("S6").Value = 0.002343443843320114
$ws.Range("T6").Value = 0.002343443843320115

# Row 7
$ws.Range("I7").Value = 0.1693441751896972
$ws.Range("J7").Value = 0.1693441751896972
$ws.Range("M7").Value = 2.124038666666666
$ws.Range("N7").Value = 6.372115999999999
$ws.Range("O7").Value = 0.4313815568786064
$ws.Range("P7").Value = 0.4313815568786064
$ws.Range("Q7").Value = 2.395732948674666
$ws.Range("R7").Value = 21.561596538072
$ws.Range("S7").Value = 0.07305195394165505
$ws.Range("T7").Value = 0.07305195394165506

# Row 8
$ws.Range("G8").Value = 5.490742
$ws.Range("H8").Value = 16.472226
$ws.Range("I8").Value = 0.8243759499123412
$ws.Range("J8").Value = 0.8243759499123412
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.731629
$ws.Range("N8").Value = 8.194887
$ws.Range("O8").Value = 0.5547800938501829
$ws.Range("P8").Value = 0.554780093850183
$ws.Range("Q8").Value = 14.998670078718
$ws.Range("R8").Value = 134.988030708462
$ws.Range("S8").Value = 0.4573473668602023
$ws.Range("T8").Value = 0.4573473668602024

# Row 9
$ws.Range("G9").Value = 5.490742
$ws.Range("H9").Value = 16.472226
$ws.Range("I9").Value = 0.8243759499123412
$ws.Range("J9").Value = 0.8243759499123412
$ws.Range("M9").Value = 0.06813733333333333
$ws.Range("O9").Value = 0.01383834927121065
$ws.Range("P9").Value = 0.01383834927121065
$ws.Range("Q9").Value = 0.3741245179013333
$ws.Range("R9").Value = 3.367120661112
$ws.Range("S9").Value = 0.01140800232567303
$ws.Range("T9").Value = 0.01140800232567303

# Row 10
$ws.Range("G10").Value = 5.490742
$ws.Range("H10").Value = 16.472226
$ws.Range("I10").Value = 0.8243759499123412
$ws.Range("J10").Value = 0.8243759499123412
$ws.Range("M10").Value = 2.124038666666666
$ws.Range("N10").Value = 6.372115999999999
$ws.Range("O10").Value = 0.4313815568786064
$ws.Range("P10").Value = 0.4313815568786064
$ws.Range("Q10").Value = 11.66254831669067
$ws.Range("R10").Value = 104.962934850216
$ws.Range("S10").Value = 0.3556205807264658
$ws.Range("T10").Value = 0.3556205807264658

